$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A61:E61").Copy()
$ws.Range("A62:E62").PasteSpecial(-4122)

$ws.Range("A62").Value = 43963
$ws.Range("B62").Value = 37351
$ws.Range("C62").Value = 1704
$ws.Range("D62").Value = 53
$ws.Range("E62").Value = 3109

$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E62"))

$ws.Range("E61").Select()
